# Split the final "rb.AddFroce(Vector3)" paragraph into two paragraphs and
# add the new "infinity" explanation text, per the commit diff.
$d = $word.ActiveDocument

# Locate the target paragraph robustly (by its current text) rather than a
# hard-coded index.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*rb.AddFroce(Vector3)*") {
        $targetPara = $candidate
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate target paragraph containing 'rb.AddFroce(Vector3)'"
}

$pr = $targetPara.Range
# Exclude the trailing paragraph mark itself so our replacement XML supplies
# the final paragraph mark (and its formatting/bookmark) explicitly.
$target = $d.Range($pr.Start, $pr.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>rb.AddFroce(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>Vector3)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>infinity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">: cuando le ponemos en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>drag</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> o en el angular </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Drag</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>infinuty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, significa que el objeto para de moverse inmediatamente.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
